$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Taniya has no pending tasks left: mark her two remaining "Working" tasks
# (row 7 and row 13 in the ToDo list) as "Done".
$ws.Range("C7").Value = "Done"
$ws.Range("C13").Value = "Done"

# Update the sheet's selection/view to reflect where the user left off.
$ws.Range("C14").Select()
